$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (pushes existing rows 3.. down by one, carrying
# formatting/comments with them - same as using Excel's Insert Row command).
$ws.Rows("3:3").Insert()

# Populate the newly inserted row with the new task.
$ws.Cells.Item(3, 1).Value = "Engine"
$ws.Cells.Item(3, 2).Value = "Refactor index based unique container.  Think about using size_t rather than unsinged int."
$ws.Cells.Item(3, 3).Value = 2

# "Basics of a physics engine" estimate changes from 21 to 35 (now on row 9
# after the insert above shifted it down from row 8).
$ws.Cells.Item(9, 3).Value = 35

# The "Scene Exporter" task row (now row 20 after the insert) is removed
# entirely.
$ws.Rows("20:20").Delete()

# Update the saved selection state.
$ws.Range("B18").Select()
